$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = [double]"24.27000000000035"
$ws.Range("H2").Value = [double]"4.000803692342906e-16"
$ws.Range("K2").Value = [double]"42.59076941581456"
$ws.Range("L2").Value = "[35.37354056381173, 49.807998267817396]"
$ws.Range("O2").Value = [double]"1.465447624197041"
$ws.Range("P2").Value = "[1.2767633807381937, 1.654131867655888]"
$ws.Range("S2").Value = [double]"62.31004276550488"
$ws.Range("T2").Value = "[57.70573334340362, 66.91435218760614]"
$ws.Range("W2").Value = [double]"18.6094294294297"
$ws.Range("X2").Value = [double]"17.88060060060086"
$ws.Range("Y2").Value = [double]"19.33825825825854"

# Row 3 updates
$ws.Range("E3").Value = [double]"25.1900000000005"
$ws.Range("G3").Value = [double]"2.220446049250313e-16"
$ws.Range("H3").Value = [double]"6.805964901916669e-16"
$ws.Range("K3").Value = [double]"33.55921784371401"
$ws.Range("L3").Value = "[24.136390407355407, 42.98204528007261]"
$ws.Range("M3").Value = [double]"1.726618847897043e-11"
$ws.Range("N3").Value = [double]"1.726618847897043e-11"
$ws.Range("O3").Value = [double]"-2.767368904063082"
$ws.Range("P3").Value = "[-3.056684744033313, -2.47805306409285]"
$ws.Range("S3").Value = [double]"55.07985016837578"
$ws.Range("T3").Value = "[49.99652482080272, 60.16317551594885]"
$ws.Range("W3").Value = [double]"11.09469469469492"
$ws.Range("X3").Value = [double]"9.934794794794994"
$ws.Range("Y3").Value = [double]"12.25459459459484"
